# CSCI6708_Case_Study_Grp1_ppt_DRAFT.pptx
# "Ethical Aspects" slide (slide 6): reword the second bullet from
# "Prevention of financial loss" to "Prevention of financial losses".
#
# The replacement re-types just the "financial loss" tail of the bullet
# (leaving the leading "Prevention of " untouched), which is what naturally
# splits the paragraph's single run into the two runs
# ("Prevention of " / "financial losses") seen in the edited file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$bulletPara = $tr.Paragraphs(3, 1)

$oldWords = "financial loss"
$newWords = "financial losses"

$fullText = $bulletPara.Text
$startPos = $fullText.IndexOf($oldWords) + 1

$target = $bulletPara.Characters($startPos, $oldWords.Length)
$target.Text = $newWords
